$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D/E target cells to Text format so numeric-looking strings
# (e.g. "1.001", "0.8088") are preserved verbatim as text, matching the
# original inline-string cell type, instead of being parsed as numbers.
$textCells = @(
    "D2", "E2", "D3", "E4", "D5", "E5", "D6", "E6", "D7", "E7",
    "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "D13",
    "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "D18", "E18",
    "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23",
    "D24", "E24", "D25", "E25", "D26", "E26", "E27", "D28", "E28", "D29",
    "E29", "D30", "E30", "D31", "E31", "E32", "D33", "E33", "D34", "E34",
    "E35", "D36", "E36", "E37", "E38", "E39", "D40", "E40", "D41", "E41",
    "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46",
    "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = '30.358.52'
$ws.Range("E2").Value = '  +1.11%  '
$ws.Range("D3").Value = '1.920.47'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '0.8088'
$ws.Range("E5").Value = '  +3.49%  '
$ws.Range("D6").Value = '244.69'
$ws.Range("E6").Value = '  +1.06%  '
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '0.3246'
$ws.Range("E8").Value = '  +2.55%  '
$ws.Range("D9").Value = '27.27'
$ws.Range("E9").Value = '  +3.51%  '
$ws.Range("D10").Value = '0.07199'
$ws.Range("E10").Value = '  +4.65%  '
$ws.Range("D11").Value = '0.7918'
$ws.Range("E11").Value = '  +6.41%  '
$ws.Range("D12").Value = '0.08096'
$ws.Range("D13").Value = '1.915.02'
$ws.Range("E13").Value = '  +1.08%  '
$ws.Range("D14").Value = '5.421'
$ws.Range("E14").Value = '  +4.20%  '
$ws.Range("D15").Value = '94.65'
$ws.Range("E15").Value = '  +1.81%  '
$ws.Range("D16").Value = '30.370.64'
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("D17").Value = '14.35'
$ws.Range("D18").Value = '6.076'
$ws.Range("E18").Value = '  +3.31%  '
$ws.Range("D19").Value = '253.62'
$ws.Range("E19").Value = '  +3.16%  '
$ws.Range("D20").Value = '0.000007845'
$ws.Range("E20").Value = '  +1.37%  '
$ws.Range("D21").Value = '2.173.44'
$ws.Range("E21").Value = '  +0.88%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.11%  '
$ws.Range("D23").Value = '8.070'
$ws.Range("E23").Value = '  +17.66%  '
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("D25").Value = '0.1630'
$ws.Range("E25").Value = '  +18.58%  '
$ws.Range("D26").Value = '9.535'
$ws.Range("E26").Value = '  +3.12%  '
$ws.Range("E27").Value = '  -0.34%  '
$ws.Range("D28").Value = '19.15'
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("D29").Value = '2.150'
$ws.Range("E29").Value = '  +5.60%  '
$ws.Range("D30").Value = '1.378'
$ws.Range("E30").Value = '  +0.76%  '
$ws.Range("D31").Value = '1.540'
$ws.Range("E31").Value = '  +1.45%  '
$ws.Range("E32").Value = '  +0.81%  '
$ws.Range("D33").Value = '4.150'
$ws.Range("E33").Value = '  +1.72%  '
$ws.Range("D34").Value = '0.05613'
$ws.Range("E34").Value = '  +0.59%  '
$ws.Range("E35").Value = '  +3.84%  '
$ws.Range("D36").Value = '0.7443'
$ws.Range("E36").Value = '  +1.44%  '
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  +1.45%  '
$ws.Range("D40").Value = '2.809'
$ws.Range("E40").Value = '  +0.74%  '
$ws.Range("D41").Value = '0.4505'
$ws.Range("E41").Value = '  +1.88%  '
$ws.Range("D42").Value = '73.65'
$ws.Range("E42").Value = '  +2.03%  '
$ws.Range("D43").Value = '5.996'
$ws.Range("E43").Value = '  -2.44%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = '1.940'
$ws.Range("E44").Value = '  +3.34%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").Value = '0.8552'
$ws.Range("E45").Value = '  +1.67%  '
$ws.Range("D46").Value = '1.000'
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").Value = '103.33'
$ws.Range("E47").Value = '  +2.90%  '
$ws.Range("D48").Value = '1.027.39'
$ws.Range("E48").Value = '  +4.14%  '
$ws.Range("D49").Value = '10.01'
$ws.Range("E49").Value = '  +2.85%  '
$ws.Range("D50").Value = '7.672'
$ws.Range("E50").Value = '  +1.83%  '
$ws.Range("D51").Value = '2.074.99'
$ws.Range("E51").Value = '  +0.94%  '

# Restore default (Normal) style on the text-formatted cells so no stray
# number-format style lingers on them once the text value is set.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
